$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New timesheet entry: row 58 ---
# Copy cell formatting from the row above (row 57) onto the new row
$ws.Range("A57:F57").Copy()
$ws.Range("A58:F58").PasteSpecial(-4122)

# Date is stored as text (like all the other date cells), keep row 57's style
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "7.4.2020"
$ws.Range("A58").NumberFormat = "General"
$ws.Range("A57").Copy()
$ws.Range("A58").PasteSpecial(-4122)

$ws.Range("B58").Value = 0.5
$ws.Range("C58").Value = 0.64583333333333337
$ws.Range("D58").Formula = "=C58-B58"

$ws.Range("E58").Value = "VGA Top"
$ws.Range("F58").Value = "Implementation"

# --- Stray formatted cell a few rows below: row 61 ---
$ws.Range("D61").NumberFormat = $ws.Range("D2").NumberFormat

# --- Selection / scrolled view moved down to the new bottom of the table ---
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("G58").Select() | Out-Null
